# Update the marksheet's "Right" (correct) count and "Total" marks values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking") - number of correct answers
$ws.Range("B11").Value = 5

# Row 12 ("Total") - total marks scored and the "scored/max" summary string
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
